# Update "想去人数" (number of people interested) counts in the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 34
$ws1.Range("F7").Value = 572
$ws1.Range("F8").Value = 8101
$ws1.Range("F9").Value = 761
$ws1.Range("F10").Value = 246
$ws1.Range("F11").Value = 1104
$ws1.Range("F12").Value = 794
$ws1.Range("F13").Value = 39
$ws1.Range("F14").Value = 32
$ws1.Range("F15").Value = 205
$ws1.Range("F16").Value = 70
$ws1.Range("F17").Value = 53
$ws1.Range("F18").Value = 210
$ws1.Range("F19").Value = 864

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 24

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 34
$ws4.Range("F8").Value = 572
$ws4.Range("F9").Value = 8102
$ws4.Range("F10").Value = 761
$ws4.Range("F11").Value = 246
$ws4.Range("F12").Value = 1104
$ws4.Range("F13").Value = 794
$ws4.Range("F14").Value = 39
$ws4.Range("F15").Value = 32
$ws4.Range("F16").Value = 205
$ws4.Range("F17").Value = 70
$ws4.Range("F18").Value = 53
$ws4.Range("F19").Value = 210
$ws4.Range("F20").Value = 864
$ws4.Range("F21").Value = 24
